# Add data for 2024-04-07
# Three new crime records were recorded for 2024 (column K / year 2024):
#   1. Near South Side  - Aggravated Assault (new crime_category for this sheet)
#   2. Jefferson Park   - Aggravated Battery
#   3. Logan Square     - Theft
# These ripple into the per-neighborhood sheets, the "By Neighborhood" summary
# sheet, and the "Citywide Totals" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Citywide Totals: bump the 2024 (K) column for the affected categories and
# the grand total.
# ---------------------------------------------------------------------------
$wsTotals = $wb.Worksheets.Item("Citywide Totals")
$wsTotals.Range("K2").Value = 29    # Aggravated Assault
$wsTotals.Range("K3").Value = 41    # Aggravated Battery
$wsTotals.Range("K10").Value = 187  # Theft
$wsTotals.Range("K11").Value = 368  # Total

# ---------------------------------------------------------------------------
# By Neighborhood: bump the 2024 (K) column for the three neighborhoods and
# the grand total row.
# ---------------------------------------------------------------------------
$wsNeigh = $wb.Worksheets.Item("By Neighborhood")
$wsNeigh.Range("K44").Value = 7    # Jefferson Park
$wsNeigh.Range("K51").Value = 3    # Logan Square
$wsNeigh.Range("K61").Value = 7    # Near South Side
$wsNeigh.Range("K97").Value = 368  # Total

# ---------------------------------------------------------------------------
# Jefferson Park: new Aggravated Battery record in 2024 (last used column is
# J because this sheet has no 2023 data).
# ---------------------------------------------------------------------------
$wsJP = $wb.Worksheets.Item("Jefferson Park")
$wsJP.Range("J3").Value = 2  # Aggravated Battery
$wsJP.Range("J6").Value = 7  # Total

# ---------------------------------------------------------------------------
# Logan Square: new Theft record in 2024.
# ---------------------------------------------------------------------------
$wsLS = $wb.Worksheets.Item("Logan Square")
$wsLS.Range("K5").Value = 2  # Theft
$wsLS.Range("K6").Value = 3  # Total

# ---------------------------------------------------------------------------
# Near South Side: this neighborhood had no prior Aggravated Assault record,
# so a brand-new row must be inserted (alphabetically, right after the
# header) and every subsequent crime_category row shifts down by one.
# Rebuild the data block explicitly so the row that disappears (the old row
# 2) and the new cell layout exactly match - a naive Rows.Insert() leaves
# stray blank styled cells across the whole row.
# ---------------------------------------------------------------------------
$wsNSS = $wb.Worksheets.Item("Near South Side")

# Drop the existing data rows (2-7) completely, keeping only the header row.
$wsNSS.Range("A2:A7").EntireRow.Delete()

# Row 2: Aggravated Assault (brand-new category for this neighborhood)
$wsNSS.Range("A2").Value = "Aggravated Assault"
$wsNSS.Range("K2").Value = 1

# Row 3: Aggravated Battery
$wsNSS.Range("A3").Value = "Aggravated Battery"
$wsNSS.Range("C3").Value = 1
$wsNSS.Range("E3").Value = 1
$wsNSS.Range("K3").Value = 1

# Row 4: Criminal Sexual Assault
$wsNSS.Range("A4").Value = "Criminal Sexual Assault"
$wsNSS.Range("J4").Value = 1

# Row 5: Homicide
$wsNSS.Range("A5").Value = "Homicide"
$wsNSS.Range("K5").Value = 1

# Row 6: Robbery
$wsNSS.Range("A6").Value = "Robbery"
$wsNSS.Range("B6").Value = 1
$wsNSS.Range("E6").Value = 3
$wsNSS.Range("F6").Value = 1
$wsNSS.Range("G6").Value = 1
$wsNSS.Range("H6").Value = 2
$wsNSS.Range("I6").Value = 3
$wsNSS.Range("J6").Value = 1
$wsNSS.Range("K6").Value = 1

# Row 7: Theft
$wsNSS.Range("A7").Value = "Theft"
$wsNSS.Range("B7").Value = 4
$wsNSS.Range("C7").Value = 3
$wsNSS.Range("D7").Value = 8
$wsNSS.Range("E7").Value = 7
$wsNSS.Range("F7").Value = 5
$wsNSS.Range("G7").Value = 3
$wsNSS.Range("I7").Value = 2
$wsNSS.Range("J7").Value = 4
$wsNSS.Range("K7").Value = 3

# Row 8: Total (2024 count increases from 6 to 7)
$wsNSS.Range("A8").Value = "Total"
$wsNSS.Range("B8").Value = 5
$wsNSS.Range("C8").Value = 4
$wsNSS.Range("D8").Value = 8
$wsNSS.Range("E8").Value = 11
$wsNSS.Range("F8").Value = 6
$wsNSS.Range("G8").Value = 4
$wsNSS.Range("H8").Value = 2
$wsNSS.Range("I8").Value = 5
$wsNSS.Range("J8").Value = 6
$wsNSS.Range("K8").Value = 7

# Re-apply the bold/centered/bordered "category" header style (same as used
# by the column headers in row 1) to the A column cells we just (re)created.
$wsNSS.Range("A1").Copy()
$wsNSS.Range("A2:A8").PasteSpecial(-4122)
